# Pareto Coefficients from WID Database
# Adds a new data row (row 4) to the estimates sheet for the
# pareto_coefficient_1990_wid estimate, fills in the correlation_direction
# helper column (J) for the existing rows, and moves the active selection
# down below the new data (matching the author's workbook state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# correlation_direction helper values for the existing rows
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 2
$ws.Range("J4").Value = 3

# New row 4: pareto_coefficient_1990_wid estimate
# (shared-string order matters for an exact match: the note text was
# authored before the estimate-name string, so write K4 first.)
$ws.Range("K4").WrapText = $true
$ws.Range("K4").Value = "As estimated by downloadParetoCoefficient.R in the taxReform1990 folder"

$ws.Range("A4").Value = "pareto_coefficient_1990_wid"
$ws.Range("B4").Value = 1.9941679999999999
$ws.Range("C4").Value = 0

$ws.Rows.Item(4).RowHeight = 30

# Move the selection to A5, matching the post-edit workbook state
$ws.Range("A5").Select()
